$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.445.94'
$ws.Range("E2").Value = '  +6.08%  '

$ws.Range("D3").Value = '1.728.64'
$ws.Range("E3").Value = '  +4.04%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9976'
$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '333.47'
$ws.Range("E5").Value = '  +5.21%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9948'
$ws.Range("E6").Value = '  -0.22%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3710'
$ws.Range("E7").Value = '  +2.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.32'
$ws.Range("E8").Value = '  +5.11%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3374'
$ws.Range("E9").Value = '  +3.18%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.202'
$ws.Range("E10").Value = '  +5.39%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07505'
$ws.Range("E11").Value = '  +6.14%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9951'
$ws.Range("E12").Value = '  -0.36%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.384'
$ws.Range("E13").Value = '  +5.45%  '

$ws.Range("E14").Value = '  +4.32%  '

$ws.Range("E15").Value = '  +5.37%  '

$ws.Range("D16").Value = '1.715.55'
$ws.Range("E16").Value = '  +3.08%  '

$ws.Range("E17").Value = '  +3.59%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06686'
$ws.Range("E18").Value = '  +0.86%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '82.99'
$ws.Range("E19").Value = '  +4.60%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9956'
$ws.Range("E20").Value = '  -0.17%  '

$ws.Range("E21").Value = '  +5.78%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.166'
$ws.Range("E22").Value = '  +4.00%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.16'
$ws.Range("E23").Value = '  +4.77%  '

$ws.Range("D24").Value = '26.399.37'
$ws.Range("E24").Value = '  +6.23%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.461'
$ws.Range("E25").Value = '  +1.31%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.507'
$ws.Range("E26").Value = '  +4.30%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.442'
$ws.Range("E27").Value = '  +17.81%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.77'
$ws.Range("E28").Value = '  +2.09%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.53'
$ws.Range("E29").Value = '  +4.87%  '

$ws.Range("D30").Value = '1.913.21'
$ws.Range("E30").Value = '  +3.53%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '131.10'
$ws.Range("E31").Value = '  +4.07%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.114'
$ws.Range("E32").Value = '  +0.70%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.075'
$ws.Range("E33").Value = '  +4.16%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08628'
$ws.Range("E34").Value = '  +2.15%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.715'
$ws.Range("E35").Value = '  +2.16%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '13.14'
$ws.Range("E36").Value = '  +6.83%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.442'
$ws.Range("E37").Value = '  +4.47%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02353'
$ws.Range("E38").Value = '  +5.18%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06331'
$ws.Range("E39").Value = '  +4.68%  '

$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2163'
$ws.Range("E40").Value = '  +4.39%  '

$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.666'
$ws.Range("E41").Value = '  +5.25%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.238'
$ws.Range("E42").Value = '  -3.18%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6243'
$ws.Range("E43").Value = '  +5.15%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.43'
$ws.Range("E44").Value = '  +13.36%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9952'
$ws.Range("E45").Value = '  -0.17%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.893'
$ws.Range("E46").Value = '  +1.90%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6032'
$ws.Range("E47").Value = '  +6.83%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '129.38'
$ws.Range("E48").Value = '  +3.42%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.057'
$ws.Range("E49").Value = '  +5.19%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07326'
$ws.Range("E50").Value = '  +4.46%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '77.59'
$ws.Range("E51").Value = '  +3.57%  '
